$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("N10").Value = "test value"
Write-Host $ws.Range("N10").Value
